$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4, shifting existing rows 4..34 down to 5..35
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new weekly price record
$ws.Range("A4").Value2 = 1
$ws.Range("B4").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value2 = "Arica y Parinacota"
$ws.Range("D4").Value2 = 45245
$ws.Range("E4").Value2 = 15
$ws.Range("F4").Value2 = 100114007
$ws.Range("G4").Value2 = "Jengibre"
$ws.Range("H4").Value2 = "Sin especificar"
$ws.Range("I4").Value2 = "Primera"
$ws.Range("J4").Value2 = 200
$ws.Range("K4").Value2 = 17000
$ws.Range("L4").Value2 = 18000
$ws.Range("M4").Value2 = 17500
$ws.Range("N4").Value2 = "`$/caja 13 kilos"
$ws.Range("O4").Value2 = "Perú"
$ws.Range("P4").Value2 = 1346
$ws.Range("Q4").Value2 = 13
$ws.Range("R4").Value2 = "Hortaliza"
